$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old data rows (A2:D11) before rewriting the restructured table
$ws.Range("A2:D11").ClearContents()

# Header: add new column E ("jenjang") with the same style as the other header cells
$ws.Range("D1").Copy()
$ws.Range("E1").PasteSpecial(-4122)
$ws.Cells.Item(1,5).Value = 'jenjang'

# Rewrite every data row (2-21) across columns A-E in the new order
$ws.Cells.Item(2,1).NumberFormat = "@"
$ws.Cells.Item(2,1).Value = '7'
$ws.Cells.Item(2,1).ClearFormats()
$ws.Cells.Item(2,2).Value = 'halaman 7'
$ws.Cells.Item(2,3).Value = 'Sekolah halaman 7'
$ws.Cells.Item(2,4).Value = 'https://drive.google.com/uc?export=download&id=1miERX46gsuFVa9PNsL56scoiK4CWqRXC'
$ws.Cells.Item(2,5).Value = 'TK'

$ws.Cells.Item(3,1).NumberFormat = "@"
$ws.Cells.Item(3,1).Value = '8'
$ws.Cells.Item(3,1).ClearFormats()
$ws.Cells.Item(3,2).Value = 'halaman 8'
$ws.Cells.Item(3,3).Value = 'Sekolah halaman 8'
$ws.Cells.Item(3,4).Value = 'https://drive.google.com/uc?export=download&id=1OtbMjfalpo7DtdUQuG31pQvwRMY2OhJO'
$ws.Cells.Item(3,5).Value = 'TK'

$ws.Cells.Item(4,1).NumberFormat = "@"
$ws.Cells.Item(4,1).Value = '3'
$ws.Cells.Item(4,1).ClearFormats()
$ws.Cells.Item(4,2).Value = 'halaman 3'
$ws.Cells.Item(4,3).Value = 'Sekolah halaman 3'
$ws.Cells.Item(4,4).Value = 'https://drive.google.com/uc?export=download&id=1WFM9dqVtzajNbJOG19HpBHmbaVsq90AA'
$ws.Cells.Item(4,5).Value = 'TK'

$ws.Cells.Item(5,1).NumberFormat = "@"
$ws.Cells.Item(5,1).Value = '5'
$ws.Cells.Item(5,1).ClearFormats()
$ws.Cells.Item(5,2).Value = 'halaman 5'
$ws.Cells.Item(5,3).Value = 'Sekolah halaman 5'
$ws.Cells.Item(5,4).Value = 'https://drive.google.com/uc?export=download&id=1-oE_Tk-l5mOIKeVMggnK2kDI6u2sM8e8'
$ws.Cells.Item(5,5).Value = 'TK'

$ws.Cells.Item(6,1).NumberFormat = "@"
$ws.Cells.Item(6,1).Value = '4'
$ws.Cells.Item(6,1).ClearFormats()
$ws.Cells.Item(6,2).Value = 'halaman 4'
$ws.Cells.Item(6,3).Value = 'Sekolah halaman 4'
$ws.Cells.Item(6,4).Value = 'https://drive.google.com/uc?export=download&id=1fuqe94jS22cc5ImQiWnKnE7VyhI0A07r'
$ws.Cells.Item(6,5).Value = 'TK'

$ws.Cells.Item(7,1).NumberFormat = "@"
$ws.Cells.Item(7,1).Value = '9'
$ws.Cells.Item(7,1).ClearFormats()
$ws.Cells.Item(7,2).Value = 'halaman 9'
$ws.Cells.Item(7,3).Value = 'Sekolah halaman 9'
$ws.Cells.Item(7,4).Value = 'https://drive.google.com/uc?export=download&id=1ExkfsNsL4bpsPn9sfYb5wJbzjv3jJ9yC'
$ws.Cells.Item(7,5).Value = 'TK'

$ws.Cells.Item(8,1).NumberFormat = "@"
$ws.Cells.Item(8,1).Value = '2'
$ws.Cells.Item(8,1).ClearFormats()
$ws.Cells.Item(8,2).Value = 'halaman 2'
$ws.Cells.Item(8,3).Value = 'Sekolah halaman 2'
$ws.Cells.Item(8,4).Value = 'https://drive.google.com/uc?export=download&id=16saOCe8EbbfEa_crnCcjdw33RR4szcZ-'
$ws.Cells.Item(8,5).Value = 'TK'

$ws.Cells.Item(9,1).NumberFormat = "@"
$ws.Cells.Item(9,1).Value = '1'
$ws.Cells.Item(9,1).ClearFormats()
$ws.Cells.Item(9,2).Value = 'halaman 1'
$ws.Cells.Item(9,3).Value = 'Sekolah halaman 1'
$ws.Cells.Item(9,4).Value = 'https://drive.google.com/uc?export=download&id=1xEeSeeEdtMcjGtSmMeCR8yuIUQ0ziuO_'
$ws.Cells.Item(9,5).Value = 'TK'

$ws.Cells.Item(10,1).NumberFormat = "@"
$ws.Cells.Item(10,1).Value = '10'
$ws.Cells.Item(10,1).ClearFormats()
$ws.Cells.Item(10,2).Value = 'halaman 10'
$ws.Cells.Item(10,3).Value = 'Sekolah halaman 10'
$ws.Cells.Item(10,4).Value = 'https://drive.google.com/uc?export=download&id=1EiOS6kSkhNZC1aGF4IlW3X36OukCZw2v'
$ws.Cells.Item(10,5).Value = 'TK'

$ws.Cells.Item(11,1).NumberFormat = "@"
$ws.Cells.Item(11,1).Value = '6'
$ws.Cells.Item(11,1).ClearFormats()
$ws.Cells.Item(11,2).Value = 'halaman 6'
$ws.Cells.Item(11,3).Value = 'Sekolah halaman 6'
$ws.Cells.Item(11,4).Value = 'https://drive.google.com/uc?export=download&id=1edRJqm7X0ciMC7sm79Tq1Z21eIMyHNJG'
$ws.Cells.Item(11,5).Value = 'TK'

$ws.Cells.Item(12,1).NumberFormat = "@"
$ws.Cells.Item(12,1).Value = '1'
$ws.Cells.Item(12,1).ClearFormats()
$ws.Cells.Item(12,2).Value = 'halaman 1'
$ws.Cells.Item(12,3).Value = 'Sekolah halaman 1'
$ws.Cells.Item(12,4).Value = 'https://drive.google.com/uc?export=download&id=1S3GHTIc8fbrVRKKrxHgtV3AvgTn6aUcA'
$ws.Cells.Item(12,5).Value = 'SD'

$ws.Cells.Item(13,1).NumberFormat = "@"
$ws.Cells.Item(13,1).Value = '4'
$ws.Cells.Item(13,1).ClearFormats()
$ws.Cells.Item(13,2).Value = 'halaman 4'
$ws.Cells.Item(13,3).Value = 'Sekolah halaman 4'
$ws.Cells.Item(13,4).Value = 'https://drive.google.com/uc?export=download&id=1y9Y6P2eo4Xs53I7yp1-FPJ_41u521Cq3'
$ws.Cells.Item(13,5).Value = 'SD'

$ws.Cells.Item(14,1).NumberFormat = "@"
$ws.Cells.Item(14,1).Value = '3'
$ws.Cells.Item(14,1).ClearFormats()
$ws.Cells.Item(14,2).Value = 'halaman 3'
$ws.Cells.Item(14,3).Value = 'Sekolah halaman 3'
$ws.Cells.Item(14,4).Value = 'https://drive.google.com/uc?export=download&id=1PAzC3xS_s7nEjtn4mdUilIXGyz3BdFUl'
$ws.Cells.Item(14,5).Value = 'SD'

$ws.Cells.Item(15,1).NumberFormat = "@"
$ws.Cells.Item(15,1).Value = '2'
$ws.Cells.Item(15,1).ClearFormats()
$ws.Cells.Item(15,2).Value = 'halaman 2'
$ws.Cells.Item(15,3).Value = 'Sekolah halaman 2'
$ws.Cells.Item(15,4).Value = 'https://drive.google.com/uc?export=download&id=1ZL7kauqEAHHD39tLnBWeHKiKE5foKx7i'
$ws.Cells.Item(15,5).Value = 'SD'

$ws.Cells.Item(16,1).NumberFormat = "@"
$ws.Cells.Item(16,1).Value = '8'
$ws.Cells.Item(16,1).ClearFormats()
$ws.Cells.Item(16,2).Value = 'halaman 8'
$ws.Cells.Item(16,3).Value = 'Sekolah halaman 8'
$ws.Cells.Item(16,4).Value = 'https://drive.google.com/uc?export=download&id=1WSYYsSDqWLLlEUJvG3cstsCvAFEXasMT'
$ws.Cells.Item(16,5).Value = 'SD'

$ws.Cells.Item(17,1).NumberFormat = "@"
$ws.Cells.Item(17,1).Value = '5'
$ws.Cells.Item(17,1).ClearFormats()
$ws.Cells.Item(17,2).Value = 'halaman 5'
$ws.Cells.Item(17,3).Value = 'Sekolah halaman 5'
$ws.Cells.Item(17,4).Value = 'https://drive.google.com/uc?export=download&id=1sRvxyZXMHYb7Sg-cal1g1EURK3eY7JFl'
$ws.Cells.Item(17,5).Value = 'SD'

$ws.Cells.Item(18,1).NumberFormat = "@"
$ws.Cells.Item(18,1).Value = '6'
$ws.Cells.Item(18,1).ClearFormats()
$ws.Cells.Item(18,2).Value = 'halaman 6'
$ws.Cells.Item(18,3).Value = 'Sekolah halaman 6'
$ws.Cells.Item(18,4).Value = 'https://drive.google.com/uc?export=download&id=1H4RGipajL8A7VeeTHm0o_UK3QG6XNIhv'
$ws.Cells.Item(18,5).Value = 'SD'

$ws.Cells.Item(19,1).NumberFormat = "@"
$ws.Cells.Item(19,1).Value = '9'
$ws.Cells.Item(19,1).ClearFormats()
$ws.Cells.Item(19,2).Value = 'halaman 9'
$ws.Cells.Item(19,3).Value = 'Sekolah halaman 9'
$ws.Cells.Item(19,4).Value = 'https://drive.google.com/uc?export=download&id=1jSTtr8O5jAynjbfiXobKoSujXf-L0IwW'
$ws.Cells.Item(19,5).Value = 'SD'

$ws.Cells.Item(20,1).NumberFormat = "@"
$ws.Cells.Item(20,1).Value = '7'
$ws.Cells.Item(20,1).ClearFormats()
$ws.Cells.Item(20,2).Value = 'halaman 7'
$ws.Cells.Item(20,3).Value = 'Sekolah halaman 7'
$ws.Cells.Item(20,4).Value = 'https://drive.google.com/uc?export=download&id=13HDcDzuPq3CW9XEzI_rXDzjmd3yHHFFj'
$ws.Cells.Item(20,5).Value = 'SD'

$ws.Cells.Item(21,1).NumberFormat = "@"
$ws.Cells.Item(21,1).Value = '10'
$ws.Cells.Item(21,1).ClearFormats()
$ws.Cells.Item(21,2).Value = 'halaman 10'
$ws.Cells.Item(21,3).Value = 'Sekolah halaman 10'
$ws.Cells.Item(21,4).Value = 'https://drive.google.com/uc?export=download&id=1F9r-UU8fvMAbgZBAs41RKipC8ra5AyCj'
$ws.Cells.Item(21,5).Value = 'SD'
